# Add "choose level" buttons row (row 37) to Sheet1, mirroring the English,
# Vietnamese "key" and Vietnamese "value" columns used by the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = "LEVEL"
$ws.Range("B37").Value = "LEVEL: "
$ws.Range("C37").Value = "ĐỘ KHÓ: "

# Match the selection left by the author after adding the new row.
$ws.Range("C37").Select()
